# Cross-browser / parallel test-run tracking: add a "TestName" column
# in front of the existing pythonCode/output columns, and label every
# test row; also drop the stray trailing duplicate "hello" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh column before A -- this slides the old A column (code
# snippets) into B and the old B column (expected/actual output) into C.
$ws.Columns("A").Insert()

# The old sheet ended with a duplicate "hello" demo row (old row 12) that
# has no corresponding output; drop it now that it lives at row 12 still.
$ws.Rows(12).Delete()

# ---- Header row ----
$ws.Range("A1").Value = "TestName"

# B1 already carries the old A1 header style (yellow fill + border) after
# the column insert; copy that same formatting onto the new C1 header so
# both header cells visually match.
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)

# The new left-most header cell gets its own look: yellow fill, no border.
$ws.Range("A1").Interior.Color = 65535

# ---- Test name labels ----
$ws.Range("A2").Value = "Print valid Hello"
$ws.Range("A3").Value = "Print invalid Hello"
$ws.Range("A4").Value = "Run Search Practice"
$ws.Range("A5").Value = "Submit Search Practice"
$ws.Range("A6").Value = "Run findMaxConsecutiveOnes Practice"
$ws.Range("A8").Value = "Run findNumbers Practice"
$ws.Range("A9").Value = "Submit findNumbers Practice"
$ws.Range("A10").Value = "Run sortedSquares Practice"
$ws.Range("A11").Value = "Submit sortedSquares Practice"
$ws.Range("A7").Value = "Submit findMaxConsecutiveOnes Practice"

# ---- Tidy up a couple of cells whose borders/shading need to line up
#      with the rest of the "plain" data cells now that the grid grew. ----
$ws.Range("C4").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("B3").PasteSpecial(-4122)

# ---- Column widths ----
$ws.Columns("A").ColumnWidth = 30.36328125
$ws.Columns("B").ColumnWidth = 195.1796875

# ---- Page setup / view ----
$ws.PageSetup.Orientation = 1
$ws.Range("A13").Select()
